# Fruta / hortaliza, semanal
# A new weekly price record is inserted as the new row 3 (date 2022-05-04 /
# serial 44685), pushing every following record down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 3, shifting rows 3-11 down to 4-12.
$ws.Rows("3").Insert()

# Populate the newly inserted row with this week's data.
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "Vega Monumental Concepción"
$ws.Range("C3").Value = "Bíobío"
$ws.Range("D3").Value = 44685
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 100112052
$ws.Range("G3").Value = "Albahaca"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 150
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 3500
$ws.Range("M3").Value = 3267
$ws.Range("N3").Value = "`$/docena de matas"
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 544
$ws.Range("Q3").Value = 6
$ws.Range("R3").Value = "Hortaliza"
